$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Cxcr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0003833333333333333
$ws.Range("N2").Value = 0.00115
$ws.Range("O2").Value = 0.005317132262509131
$ws.Range("P2").Value = 0.005317132262509131
$ws.Range("Q2").Value = 0.06513433518333334
$ws.Range("R2").Value = 0.58620901665
$ws.Range("S2").Value = 0.002361436768429497
$ws.Range("T2").Value = 0.002361436768429497

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Cxcr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.07171066666666667
$ws.Range("N3").Value = 0.215132
$ws.Range("O3").Value = 0.9946828677374909
$ws.Range("P3").Value = 0.9946828677374908
$ws.Range("Q3").Value = 12.18476504057467
$ws.Range("R3").Value = 109.662885365172
$ws.Range("S3").Value = 0.4417570564050214
$ws.Range("T3").Value = 0.4417570564050214

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Cxcr1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 68.382243
$ws.Range("H4").Value = 205.146729
$ws.Range("I4").Value = 0.1787346690539575
$ws.Range("J4").Value = 0.1787346690539575
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0003833333333333333
$ws.Range("N4").Value = 0.00115
$ws.Range("O4").Value = 0.005317132262509131
$ws.Range("P4").Value = 0.005317132262509131
$ws.Range("Q4").Value = 0.02621319315
$ws.Range("R4").Value = 0.23591873835
$ws.Range("S4").Value = 0.0009503558752556899
$ws.Range("T4").Value = 0.0009503558752556899

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Cxcr1"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 68.382243
$ws.Range("H5").Value = 205.146729
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.07171066666666667
$ws.Range("N5").Value = 0.215132
$ws.Range("O5").Value = 0.9946828677374909
$ws.Range("P5").Value = 0.9946828677374908
$ws.Range("Q5").Value = 4.903736233692
$ws.Range("R5").Value = 44.133626103228
$ws.Range("S5").Value = 0.1777843131787018
$ws.Range("T5").Value = 0.1777843131787018

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Cxcr1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 53.27463399999999
$ws.Range("H6").Value = 159.823902
$ws.Range("I6").Value = 0.1392470275793777
$ws.Range("J6").Value = 0.1392470275793778
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.0003833333333333333
$ws.Range("N6").Value = 0.00115
$ws.Range("O6").Value = 0.005317132262509131
$ws.Range("P6").Value = 0.005317132262509131
$ws.Range("Q6").Value = 0.02042194303333333
$ws.Range("R6").Value = 0.1837974873
$ws.Range("S6").Value = 0.0007403948628008081
$ws.Range("T6").Value = 0.0007403948628008082

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Cxcr1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 53.27463399999999
$ws.Range("H7").Value = 159.823902
$ws.Range("I7").Value = 0.1392470275793777
$ws.Range("J7").Value = 0.1392470275793778
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.07171066666666667
$ws.Range("N7").Value = 0.215132
$ws.Range("O7").Value = 0.9946828677374909
$ws.Range("P7").Value = 0.9946828677374908
$ws.Range("Q7").Value = 3.820359520562667
$ws.Range("R7").Value = 34.383235685064
$ws.Range("S7").Value = 0.1385066327165769
$ws.Range("T7").Value = 0.138506632716577

$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Cxcr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 91.01828266666666
$ws.Range("H8").Value = 273.054848
$ws.Range("I8").Value = 0.2378998101932138
$ws.Range("J8").Value = 0.2378998101932138
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.0003833333333333333
$ws.Range("N8").Value = 0.00115
$ws.Range("O8").Value = 0.005317132262509131
$ws.Range("P8").Value = 0.005317132262509131
$ws.Range("Q8").Value = 0.03489034168888889
$ws.Range("R8").Value = 0.3140130752
$ws.Range("S8").Value = 0.001264944756023136
$ws.Range("T8").Value = 0.001264944756023136

$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Cxcr1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 91.01828266666666
$ws.Range("H9").Value = 273.054848
$ws.Range("I9").Value = 0.2378998101932138
$ws.Range("J9").Value = 0.2378998101932138
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.07171066666666667
$ws.Range("N9").Value = 0.215132
$ws.Range("O9").Value = 0.9946828677374909
$ws.Range("P9").Value = 0.9946828677374908
$ws.Range("Q9").Value = 6.526981728881778
$ws.Range("R9").Value = 58.742835559936
$ws.Range("S9").Value = 0.2366348654371907
$ws.Range("T9").Value = 0.2366348654371907

